# Updates the cryptocurrency price/volume table (columns D, E) with the
# latest scraped values, mirroring the GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.672.23"
$ws.Range("E2").Value = "  -1.44%  "
$ws.Range("D3").Value = "1.594.20"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.09"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.65"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0835"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "1.818.98"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "1.601.61"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "26.648.38"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.34"
$ws.Range("E19").Value = "  -2.59%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.58"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -3.61%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.31"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.668"
$ws.Range("E33").Value = "  -4.50%  "
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").Value = "1.300.14"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("E37").Value = "  -4.72%  "
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.836"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.793"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.53"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "1.730.84"
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.907"
$ws.Range("E46").Value = "  +6.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.75"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.47"
$ws.Range("E51").Value = "  -1.36%  "
